# "Updated project journal time"
#
# Jesse logged a new entry in his personal time-tracking tab: on 10/30/2017
# he spent 60 minutes "Added new header files. Implemented part of Database.h".
# That entry lands in the first open row (row 5) below the existing log row
# (row 4) on the "Jesse" sheet. The sheet's running total (C2, =SUM(B4:B200))
# and the roll-up on the "Main" sheet (B2, =(Jesse!C2)/60) both pick the new
# entry up automatically via recalculation.

$wb = $excel.ActiveWorkbook
$jesse = $wb.Worksheets.Item("Jesse")
$main = $wb.Worksheets.Item("Main")

$jesse.Range("A5").Value = 43038   # Date -> 10/30/2017
$jesse.Range("B5").Value = 60      # Time Spent (minutes)
$jesse.Range("C5").Value = "Added new header files. Implemented part of Database.h"

# The description wraps onto two lines at this column width, so the row grows
# to fit it.
$jesse.Rows.Item(5).RowHeight = 28.5

# Leave the cursor where Jesse last clicked while filling in the new row...
$jesse.Range("G13:G14").Select() | Out-Null

# ...then hop back to the "Main" overview tab, which is what's on screen
# when the workbook is saved.
$main.Activate()
